# "SOLVED" version of the GradeBook activity spreadsheet.
# Fills in the Final Grade (F), Pass/Fail (G), and Letter Grade (H)
# columns with the worked-out formulas, then leaves the selection on G10
# (matching the author's saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the "master" (non-shared) formulas...
$ws.Range("F2").Formula = "=ROUND(AVERAGE(B2:E2),0)"
$ws.Range("G2").Formula = '=IF(F2>60,"PASS","FAIL")'
$ws.Range("H2").Formula = '=IF(F2>=90,"A",IF(F2>=80,"B",IF(F2>=70,"C",IF(F2>=60,"D","F"))))'

# ...while rows 3-25 share the fill-down formula relative to row 3.
$ws.Range("F3:F25").Formula = "=ROUND(AVERAGE(B3:E3),0)"
$ws.Range("G3:G25").Formula = '=IF(F3>60,"PASS","FAIL")'
$ws.Range("H3:H25").Formula = '=IF(F3>=90,"A",IF(F3>=80,"B",IF(F3>=70,"C",IF(F3>=60,"D","F"))))'

# Move the active selection to G10, matching the saved workbook state.
$ws.Range("G10").Select() | Out-Null
